$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 800
    5  = 600
    6  = 650
    8  = 400
    9  = 700
    10 = 400
    12 = 600
    13 = 600
    14 = 800
    15 = 650
    18 = 700
    19 = 400
    20 = 500
    21 = 700
    22 = 800
    23 = 800
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}
